# "New data for 100" — updates the Trial values in column D (the "100"
# input-size column) of the Time Analysis sheet, and moves the sheet
# selection/scroll position to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the trial counts for input size 100 (see D5 = 100).
# New measurements replaced the previous values for these rows.
$ws.Range("D9").Value  = 0
$ws.Range("D11").Value = 2
$ws.Range("D12").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("D16").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("D21").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 11
$ws.Range("D27").Value = 9
$ws.Range("D28").Value = 5
$ws.Range("D29").Value = 7
$ws.Range("D30").Value = 3
$ws.Range("D31").Value = 1
$ws.Range("D36").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("D43").Value = 1
$ws.Range("D44").Value = 1
$ws.Range("D46").Value = 1
$ws.Range("D48").Value = 1
$ws.Range("D49").Value = 1
$ws.Range("D53").Value = 1
$ws.Range("D54").Value = 1
$ws.Range("D58").Value = 0
$ws.Range("D61").Value = 1
$ws.Range("D64").Value = 0

# Reflect the author's final scroll position/selection on the sheet.
$ws.Activate()
$ws.Range("D66").Select()
